# Atualização de bases das ligas, do dia: 11-04-2024 às 00:31
#
# The fixtures in rows 8-10 (match ids 6,7,8) had their HomeTeam/AwayTeam
# and odds data rotated by one row; rows 93-94 (match ids 91,92) and
# rows 157-158 (match ids 155,156) had their data swapped. The "id"
# column (A) stays put - only columns B:AC move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 8, 9, 10: cyclic rotation ---------------------------------
# after(8) = before(9); after(9) = before(10); after(10) = before(8)
$row8 = $ws.Range("B8:AC8")
$row9 = $ws.Range("B9:AC9")
$row10 = $ws.Range("B10:AC10")

$v8 = $row8.Value()
$v9 = $row9.Value()
$v10 = $row10.Value()

$row8.Value = $v9
$row9.Value = $v10
$row10.Value = $v8

# --- Rows 93, 94: swap ------------------------------------------------
$row93 = $ws.Range("B93:AC93")
$row94 = $ws.Range("B94:AC94")

$v93 = $row93.Value()
$v94 = $row94.Value()

$row93.Value = $v94
$row94.Value = $v93

# --- Rows 157, 158: swap ----------------------------------------------
$row157 = $ws.Range("B157:AC157")
$row158 = $ws.Range("B158:AC158")

$v157 = $row157.Value()
$v158 = $row158.Value()

$row157.Value = $v158
$row158.Value = $v157
